$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.439655512935477
$ws.Range("C2").Value = 4.910019341749957
$ws.Range("E2").Value = 23.46257738199308
$ws.Range("F2").Value = 40.68532829646367
$ws.Range("G2").Value = 3.638907270510543
$ws.Range("I2").Value = 20.55930814898493
$ws.Range("J2").Value = 8.062935027975461
$ws.Range("K2").Value = 9.258557754030093
$ws.Range("O2").Value = 21.8838791473291

$ws.Range("B3").Value = 9.081829350803954
$ws.Range("C3").Value = 4.664613401182223
$ws.Range("E3").Value = 23.09693550480246
$ws.Range("F3").Value = 40.54071422970096
$ws.Range("G3").Value = 3.640572480825728
$ws.Range("I3").Value = 20.68878718183965
$ws.Range("J3").Value = 8.088791244318241
$ws.Range("K3").Value = 9.009781126830694
$ws.Range("O3").Value = 22.01107926560912

$ws.Range("B4").Value = 8.855420559057553
$ws.Range("C4").Value = 4.506290014138854
$ws.Range("E4").Value = 22.87475314712606
$ws.Range("F4").Value = 40.46349433230809
$ws.Range("G4").Value = 3.641647646974465
$ws.Range("I4").Value = 20.77261681874335
$ws.Range("J4").Value = 8.10555254767069
$ws.Range("K4").Value = 8.854211390313131
$ws.Range("O4").Value = 22.094291075222

$ws.Range("B5").Value = 8.761603321391672
$ws.Range("C5").Value = 4.439893593130865
$ws.Range("E5").Value = 22.78490904697793
$ws.Range("F5").Value = 40.43495796752435
$ws.Range("G5").Value = 3.642099084810587
$ws.Range("I5").Value = 20.8078681155033
$ws.Range("J5").Value = 8.112606100719018
$ws.Range("K5").Value = 8.790194002216152
$ws.Range("O5").Value = 22.12948441686786

$ws.Range("B6").Value = 8.745935320322303
$ws.Range("C6").Value = 4.428756507259179
$ws.Range("E6").Value = 22.77003588431531
$ws.Range("F6").Value = 40.43039711043534
$ws.Range("G6").Value = 3.642174850174194
$ws.Range("I6").Value = 20.81378744654505
$ws.Range("J6").Value = 8.113790833177605
$ws.Range("K6").Value = 8.779529093991986
$ws.Range("O6").Value = 22.13540575243602

$ws.Range("B7").Value = 8.854161413980329
$ws.Range("C7").Value = 4.505402110533137
$ws.Range("E7").Value = 22.87353850743738
$ws.Range("F7").Value = 40.46309758713959
$ws.Range("G7").Value = 3.641653681320425
$ws.Range("I7").Value = 20.77308781455507
$ws.Range("J7").Value = 8.105646769986077
$ws.Range("K7").Value = 8.85335042600671
$ws.Range("O7").Value = 22.09476050821957

$ws.Range("B8").Value = 9.317742521162153
$ws.Range("C8").Value = 4.827013384287024
$ws.Range("E8").Value = 23.33608968966927
$ws.Range("F8").Value = 40.63307707708967
$ws.Range("G8").Value = 3.639470518251131
$ws.Range("I8").Value = 20.60305492494781
$ws.Range("J8").Value = 8.071666785345871
$ws.Range("K8").Value = 9.173415737357018
$ws.Range("O8").Value = 21.92667634118557

$ws.Range("B9").Value = 10.16858339745011
$ws.Range("C9").Value = 5.395516673681337
$ws.Range("E9").Value = 24.25662527660384
$ws.Range("F9").Value = 41.05712989660161
$ws.Range("G9").Value = 3.635605731454915
$ws.Range("I9").Value = 20.30390161862786
$ws.Range("J9").Value = 8.012034112518556
$ws.Range("K9").Value = 9.77511865052767
$ws.Range("O9").Value = 21.63766776754769

$ws.Range("B10").Value = 10.75234500536061
$ws.Range("C10").Value = 5.773691579611321
$ws.Range("E10").Value = 24.93457962370504
$ws.Range("F10").Value = 41.42224654864729
$ws.Range("G10").Value = 3.63301741243591
$ws.Range("I10").Value = 20.10491679464077
$ws.Range("J10").Value = 7.972457082477563
$ws.Range("K10").Value = 10.19700045143658
$ws.Range("O10").Value = 21.45016398502469

$ws.Range("B11").Value = 11.00794180813366
$ws.Range("C11").Value = 5.936921419160597
$ws.Range("E11").Value = 25.24199954107389
$ws.Range("F11").Value = 41.59953047231339
$ws.Range("G11").Value = 3.631893875934849
$ws.Range("I11").Value = 20.01889142133278
$ws.Range("J11").Value = 7.955365014122992
$ws.Range("K11").Value = 10.38373424028216
$ws.Range("O11").Value = 21.37027409908191

$ws.Range("B12").Value = 11.10322910367878
$ws.Range("C12").Value = 5.997451986701032
$ws.Range("E12").Value = 25.35816096357625
$ws.Range("F12").Value = 41.66822939411228
$ws.Range("G12").Value = 3.631476129585341
$ws.Range("I12").Value = 19.98696083657498
$ws.Range("J12").Value = 7.949023275919894
$ws.Range("K12").Value = 10.45364076970644
$ws.Range("O12").Value = 21.34080121873465

$ws.Range("B13").Value = 11.08277507325545
$ws.Range("C13").Value = 5.9844728246625
$ws.Range("E13").Value = 25.33315643475486
$ws.Range("F13").Value = 41.65336496400325
$ws.Range("G13").Value = 3.631565756330399
$ws.Range("I13").Value = 19.99380896959135
$ws.Range("J13").Value = 7.950383278744022
$ws.Range("K13").Value = 10.43862187297733
$ws.Range("O13").Value = 21.34711401738135

$ws.Range("B14").Value = 11.0158116114555
$ws.Range("C14").Value = 5.941927053551114
$ws.Range("E14").Value = 25.25156187516953
$ws.Range("F14").Value = 41.6051512642664
$ws.Range("G14").Value = 3.631859353340928
$ws.Range("I14").Value = 20.01625154358472
$ws.Range("J14").Value = 7.954840659477071
$ws.Range("K14").Value = 10.38950196134404
$ws.Range("O14").Value = 21.36783370823199

$ws.Range("B15").Value = 10.97459709530476
$ws.Range("C15").Value = 5.915699322626187
$ws.Range("E15").Value = 25.20154685353043
$ws.Range("F15").Value = 41.57582148090941
$ws.Range("G15").Value = 3.632040193245685
$ws.Range("I15").Value = 20.03008228880871
$ws.Range("J15").Value = 7.957587936145075
$ws.Range("K15").Value = 10.35930800374317
$ws.Range("O15").Value = 21.38062672339782

$ws.Range("B16").Value = 10.73543458259067
$ws.Range("C16").Value = 5.762845605858924
$ws.Range("E16").Value = 24.914459680526
$ws.Range("F16").Value = 41.41088213453592
$ws.Range("G16").Value = 3.633091919542875
$ws.Range("I16").Value = 20.11062910798976
$ws.Range("J16").Value = 7.973592397104552
$ws.Range("K16").Value = 10.18468739411761
$ws.Range("O16").Value = 21.45549395852364

$ws.Range("B17").Value = 10.58611274487536
$ws.Range("C17").Value = 5.666808121483443
$ws.Range("E17").Value = 24.73800910710263
$ws.Range("F17").Value = 41.31253317840422
$ws.Range("G17").Value = 3.633750897792717
$ws.Range("I17").Value = 20.16119232038041
$ws.Range("J17").Value = 7.983643808664814
$ws.Range("K17").Value = 10.07619128132198
$ws.Range("O17").Value = 21.50280887799917

$ws.Range("B18").Value = 10.49929286551192
$ws.Range("C18").Value = 5.610741904993526
$ws.Range("E18").Value = 24.63643306152781
$ws.Range("F18").Value = 41.25702159914253
$ws.Range("G18").Value = 3.634135000529625
$ws.Range("I18").Value = 20.19069796853252
$ws.Range("J18").Value = 7.989510956528696
$ws.Range("K18").Value = 10.01330263613941
$ws.Range("O18").Value = 21.53053191548566

$ws.Range("B19").Value = 10.46973902432339
$ws.Range("C19").Value = 5.591617153315385
$ws.Range("E19").Value = 24.60202974982786
$ws.Range("F19").Value = 41.23840900861234
$ws.Range("G19").Value = 3.634265924084665
$ws.Range("I19").Value = 20.20076075622218
$ws.Range("J19").Value = 7.99151222914268
$ws.Range("K19").Value = 9.99192832930669
$ws.Range("O19").Value = 21.54000577122289

$ws.Range("B20").Value = 10.60210549949337
$ws.Range("C20").Value = 5.677117289796239
$ws.Range("E20").Value = 24.75680231876703
$ws.Range("F20").Value = 41.32289359367001
$ws.Range("G20").Value = 3.633680223379052
$ws.Range("I20").Value = 20.15576600488876
$ws.Range("J20").Value = 7.982564937702488
$ws.Range("K20").Value = 10.08779148154919
$ws.Range("O20").Value = 21.49771945475269

$ws.Range("B21").Value = 11.03552166527853
$ws.Range("C21").Value = 5.954458648135329
$ws.Range("E21").Value = 25.27553586807952
$ws.Range("F21").Value = 41.61927067836469
$ws.Range("G21").Value = 3.631772907772324
$ws.Range("I21").Value = 20.00964210717042
$ws.Range("J21").Value = 7.953527876238938
$ws.Range("K21").Value = 10.40395198377612
$ws.Range("O21").Value = 21.3617266539842

$ws.Range("B22").Value = 11.31000556950657
$ws.Range("C22").Value = 6.128245687494218
$ws.Range("E22").Value = 25.61304219695458
$ws.Range("F22").Value = 41.82207179539814
$ws.Range("G22").Value = 3.630571303626058
$ws.Range("I22").Value = 19.91790298818397
$ws.Range("J22").Value = 7.935311845353131
$ws.Range("K22").Value = 10.60586201238953
$ws.Range("O22").Value = 21.27739321646717

$ws.Range("B23").Value = 11.16433214961902
$ws.Range("C23").Value = 6.036180061013327
$ws.Range("E23").Value = 25.43308278154129
$ws.Range("F23").Value = 41.71301550642964
$ws.Range("G23").Value = 3.631208523333566
$ws.Range("I23").Value = 19.96652198457694
$ws.Range("J23").Value = 7.944964564605028
$ws.Range("K23").Value = 10.49854910538565
$ws.Range("O23").Value = 21.32198685905582

$ws.Range("B24").Value = 10.59487819969961
$ws.Range("C24").Value = 5.672459170531802
$ws.Range("E24").Value = 24.7483063112772
$ws.Range("F24").Value = 41.31820643646014
$ws.Range("G24").Value = 3.633712158933268
$ws.Range("I24").Value = 20.1582178833315
$ws.Range("J24").Value = 7.98305241968851
$ws.Range("K24").Value = 10.08254862358803
$ws.Range("O24").Value = 21.50001876005211

$ws.Range("B25").Value = 9.945295756806507
$ws.Range("C25").Value = 5.248564593333998
$ws.Range("E25").Value = 24.00683927787074
$ws.Range("F25").Value = 40.93286473376473
$ws.Range("G25").Value = 3.636606961733262
$ws.Range("I25").Value = 20.38117063762682
$ws.Range("J25").Value = 8.027420209491432
$ws.Range("K25").Value = 9.615593956018524
$ws.Range("O25").Value = 21.71149768604512
